# Scheduled data refresh: update market-price derived columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3863.7778
$ws.Range("I62").Value = 5718.3335
$ws.Range("J62").Value = 2380.1333
$ws.Range("K62").Value = 5718.3335
$ws.Range("L62").Value = 2380.1333
$ws.Range("M62").Value = -5094.3335
$ws.Range("N62").Value = -3628.1333

$ws.Range("H65").Value = 3863.7778
$ws.Range("I65").Value = 5718.3335
$ws.Range("J65").Value = 2380.1333
$ws.Range("K65").Value = 28591.6675
$ws.Range("L65").Value = 11900.6665
$ws.Range("M65").Value = -25471.6675
$ws.Range("N65").Value = -18140.6665

$ws.Range("H107").Value = 659.04
$ws.Range("I107").Value = 691
$ws.Range("J107").Value = 637.73334
$ws.Range("K107").Value = 691
$ws.Range("L107").Value = 637.73334
$ws.Range("M107").Value = 1229
$ws.Range("N107").Value = -4477.73334

$ws.Range("H116").Value = 2969.5652
$ws.Range("I116").Value = 2664.2856
$ws.Range("K116").Value = 2664.2856
$ws.Range("M116").Value = 777.7143999999998

$ws.Range("H138").Value = 1843.2
$ws.Range("I138").Value = 1153.0869
$ws.Range("J138").Value = 2049.3376
$ws.Range("K138").Value = 3459.2607
$ws.Range("L138").Value = 6148.0128
$ws.Range("M138").Value = 1680.7393
$ws.Range("N138").Value = -16428.0128

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1406.28
$ws.Range("I2").Value = 1445.8695
$ws.Range("K2").Value = 1445.8695
$ws.Range("M2").Value = -1332.8695

$ws.Range("H61").Value = 1760.5927
$ws.Range("I61").Value = 1420.5238
$ws.Range("J61").Value = 2950.8333
$ws.Range("K61").Value = 1420.5238
$ws.Range("L61").Value = 2950.8333
$ws.Range("M61").Value = -1208.5238
$ws.Range("N61").Value = -3374.8333

$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

$ws.Range("H110").Value = 1228.5682
$ws.Range("I110").Value = 1157.5264
$ws.Range("J110").Value = 1678.5
$ws.Range("K110").Value = 1157.5264
$ws.Range("L110").Value = 1678.5
$ws.Range("M110").Value = 887.4736
$ws.Range("N110").Value = -5768.5

$ws.Range("H116").Value = 1406.28
$ws.Range("I116").Value = 1445.8695
$ws.Range("K116").Value = 1445.8695
$ws.Range("M116").Value = 848.1305

$ws.Range("H136").Value = 1760.5927
$ws.Range("I136").Value = 1420.5238
$ws.Range("J136").Value = 2950.8333
$ws.Range("K136").Value = 4261.5714
$ws.Range("L136").Value = 8852.499899999999
$ws.Range("M136").Value = -1711.5714
$ws.Range("N136").Value = -13952.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1406.28
$ws.Range("I3").Value = 1445.8695
$ws.Range("K3").Value = 1445.8695
$ws.Range("M3").Value = -1331.8695

$ws.Range("H50").Value = 19955
$ws.Range("J50").Value = 19955
$ws.Range("L50").Value = 19955
$ws.Range("N50").Value = -21103

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H134").Value = 2144.3723
$ws.Range("I134").Value = 1373.5593
$ws.Range("J134").Value = 3443.743
$ws.Range("K134").Value = 4120.6779
$ws.Range("L134").Value = 10331.229
$ws.Range("M134").Value = -1585.6779
$ws.Range("N134").Value = -15401.229

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2731.276
$ws.Range("I105").Value = 2627.5908
$ws.Range("J105").Value = 3057.1428
$ws.Range("K105").Value = 2627.5908
$ws.Range("L105").Value = 3057.1428
$ws.Range("M105").Value = -880.5907999999999
$ws.Range("N105").Value = -6551.1428

$ws.Range("H107").Value = 800.7941
$ws.Range("I107").Value = 555.1818
$ws.Range("J107").Value = 1251.0834
$ws.Range("K107").Value = 555.1818
$ws.Range("L107").Value = 1251.0834
$ws.Range("M107").Value = 1364.8182
$ws.Range("N107").Value = -5091.0834

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3560.125
$ws.Range("I5").Value = 3610.4517
$ws.Range("K5").Value = 10831.3551
$ws.Range("M5").Value = -10719.3551

$ws.Range("H113").Value = 5774.45
$ws.Range("I113").Value = 7687.5
$ws.Range("J113").Value = 1310.6666
$ws.Range("K113").Value = 23062.5
$ws.Range("L113").Value = 3931.9998
$ws.Range("M113").Value = -20892.5
$ws.Range("N113").Value = -8271.9998

$ws.Range("H131").Value = 865.74
$ws.Range("I131").Value = 482.25
$ws.Range("J131").Value = 881.71875
$ws.Range("K131").Value = 1446.75
$ws.Range("L131").Value = 2645.15625
$ws.Range("M131").Value = 3593.25
$ws.Range("N131").Value = -12725.15625

$ws.Range("H135").Value = 3560.125
$ws.Range("I135").Value = 3610.4517
$ws.Range("K135").Value = 32494.0653
$ws.Range("M135").Value = -29959.0653

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2206.5305
$ws.Range("I132").Value = 1566.909
$ws.Range("J132").Value = 3525.75
$ws.Range("K132").Value = 4700.727000000001
$ws.Range("L132").Value = 10577.25
$ws.Range("M132").Value = -2170.727000000001
$ws.Range("N132").Value = -15637.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6411590
$ws.Range("I82").Value = 1304.4445
$ws.Range("J82").Value = 20834732
$ws.Range("K82").Value = 1304.4445
$ws.Range("L82").Value = 20834732
$ws.Range("M82").Value = -943.4445000000001
$ws.Range("N82").Value = -20835454

$ws.Range("H85").Value = 6411590
$ws.Range("I85").Value = 1304.4445
$ws.Range("J85").Value = 20834732
$ws.Range("K85").Value = 1304.4445
$ws.Range("L85").Value = 20834732
$ws.Range("M85").Value = -56.44450000000006
$ws.Range("N85").Value = -20837228

$ws.Range("H93").Value = 2053.4666
$ws.Range("I93").Value = 2599.3333
$ws.Range("J93").Value = 1917
$ws.Range("K93").Value = 2599.3333
$ws.Range("L93").Value = 1917
$ws.Range("M93").Value = -1351.3333
$ws.Range("N93").Value = -4413

$ws.Range("H122").Value = 2229.9
$ws.Range("I122").Value = 2271.0588
$ws.Range("K122").Value = 6813.176399999999
$ws.Range("M122").Value = -4363.176399999999

$ws.Range("H132").Value = 2544.1428
$ws.Range("I132").Value = 1765.4073
$ws.Range("J132").Value = 3945.8667
$ws.Range("K132").Value = 5296.2219
$ws.Range("L132").Value = 11837.6001
$ws.Range("M132").Value = -2766.2219
$ws.Range("N132").Value = -16897.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3413.2727
$ws.Range("I81").Value = 872.6667
$ws.Range("J81").Value = 8857.429
$ws.Range("K81").Value = 1745.3334
$ws.Range("L81").Value = 17714.858
$ws.Range("M81").Value = -684.3334
$ws.Range("N81").Value = -19836.858

$ws.Range("H84").Value = 3413.2727
$ws.Range("I84").Value = 872.6667
$ws.Range("J84").Value = 8857.429
$ws.Range("K84").Value = 8726.666999999999
$ws.Range("L84").Value = 88574.29000000001
$ws.Range("M84").Value = -3422.666999999999
$ws.Range("N84").Value = -99182.29000000001

$ws.Range("H92").Value = 33030
$ws.Range("J92").Value = 33030
$ws.Range("L92").Value = 33030
$ws.Range("N92").Value = -38022

$ws.Range("H107").Value = 6250777
$ws.Range("I107").Value = 302.66666
$ws.Range("K107").Value = 907.9999799999999
$ws.Range("M107").Value = 1012.00002

$ws.Range("H113").Value = 359.16666
$ws.Range("I113").Value = 382.72726
$ws.Range("J113").Value = 100
$ws.Range("K113").Value = 1148.18178
$ws.Range("L113").Value = 300
$ws.Range("M113").Value = 1021.81822
$ws.Range("N113").Value = -4640
